$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target layout for rows 2-68 (date, col B, col C, col D, col E).
# Rows 2-49 already existed but are being re-ordered (each year's Oct/Nov/Dec
# rows move to the front of that year's block); rows 50-68 are brand-new data
# appended through 2023-07.
$rows = @(
    @("2018-10", 106.6, 109.8, 106.6, 109.2),
    @("2018-11", 104.7, 111, 104.8, 108.4),
    @("2018-12", 98.6, 104.7, 96.1, 105.4),
    @("2018-01", 120.5031, 114.3415, 115.2722, 103.9098),
    @("2018-02", 114.9, 110, 111.5, 106.3),
    @("2018-03", 112.7, 108.3, 109.1, 107.4),
    @("2018-04", 113.4, 106.9, 111.4, 105.9),
    @("2018-05", 112.9, 109.2, 115.2, 107.1),
    @("2018-06", 114, 112, 115.7, 109.3),
    @("2018-07", 109.5, 111.8, 112.7, 111.8),
    @("2018-08", 107, 111, 109.5, 110.8),
    @("2018-09", 106.5, 110, 107, 109.7),
    @("2019-10", 94.2, 102.1, 93.6, 91.9),
    @("2019-11", 95.7, 101.4, 94.6, 90.2),
    @("2019-12", 99.4, 103.7, 99.8, 90.1),
    @("2019-01", 99.1, 102.7, 96.3, 100.3),
    @("2019-02", 100, 103.3, 97.4, 99),
    @("2019-03", 99.8, 102.8, 98.1, 98.2),
    @("2019-04", 103, 103.1, 102.5, 99.5),
    @("2019-05", 103.7, 103.5, 101.8, 99.5),
    @("2019-06", 100.9, 103.4, 99.1, 99),
    @("2019-07", 101.8, 105.5, 99.3, 96.8),
    @("2019-08", 98.3, 104.8, 96.5, 95),
    @("2019-09", 94.9, 102.1, 93.8, 92.7),
    @("2020-10", 100.7, 100, 100.8, 93.2),
    @("2020-11", 101.2, 100.7, 102.8, 94.8),
    @("2020-12", 103.3, 102.1, 105.9, 99),
    @("2020-01", 99.8, 104.2, 101.3, 92),
    @("2020-02", 96.5, 103.4, 99.6, 93.3),
    @("2020-03", 95, 101.7, 96.2, 92.5),
    @("2020-04", 92.6, 100.8, 92.3, 91.6),
    @("2020-05", 92.9, 99.3, 91.4, 92.2),
    @("2020-06", 96.3, 99.3, 94.7, 92.3),
    @("2020-07", 95.7, 97.8, 95.6, 90.9),
    @("2020-08", 98.6, 97.9, 97.9, 90.4),
    @("2020-09", 100.9, 99.7, 100.6, 91.9),
    @("2021-10", 140.3, 135, 138.5, 162.6),
    @("2021-11", 134.5, 130.3, 128.1, 165.9),
    @("2021-12", 127.8, 121.2, 118.5, 149),
    @("2021-01", 107.6, 110.8, 110.7, 104.4),
    @("2021-02", 114.9, 114.2, 114.5, 106.9),
    @("2021-03", 123.4, 118, 122, 112.7),
    @("2021-04", 134.2, 122.4, 130.9, 112.1),
    @("2021-05", 141, 132.1, 139.7, 112.2),
    @("2021-06", 134.1, 134.6, 135.8, 115.4),
    @("2021-07", 135.2, 135.2, 133.3, 121.6),
    @("2021-08", 133.4, 132.2, 134.6, 129),
    @("2021-09", 134.8, 132.2, 134.8, 139),
    @("2022-10", 81, 83.2, 78.3, 80),
    @("2022-11", 81.6, 83.3, 81.3, 79.1),
    @("2022-12", 84.3, 88.2, 85.4, 85.9),
    @("2022-01", 123.2, 115.7, 111.7, 139.5),
    @("2022-02", 119.9, 113.1, 110.5, 132.8),
    @("2022-03", 113.7, 110.8, 107.4, 128.5),
    @("2022-04", 107.6, 111.3, 103.8, 126.8),
    @("2022-05", 100.1, 104.8, 96.7, 126.2),
    @("2022-06", 100.7, 97.7, 94.1, 120.3),
    @("2022-07", 92.5, 91.5, 88.8, 111.8),
    @("2022-08", 87.4, 88.4, 83.3, 101.9),
    @("2022-09", 85.5, 85.6, 80.6, 92.4),
    @("2023-01", 86.6, 89.8, 88.6, 88.4),
    @("2023-02", 87, 89.2, 89.4, 90.6),
    @("2023-03", 87.3, 87.7, 89.2, 89.7),
    @("2023-04", 85.9, 84.9, 86.4, 88.7),
    @("2023-05", 82.4, 80.3, 83.2, 86.2),
    @("2023-06", 83.1, 81.4, 84, 86.7),
    @("2023-07", 89.5, 85.7, 89.6, 88)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# Rows 50-68 are brand new; copy the date-column formatting (bold, centred,
# thin border) from an existing formatted cell so the new cells match style "s=1".
$ws.Range("A49").Copy() | Out-Null
$ws.Range("A50:A68").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false